# Add a new weekly ranking sheet "2025-12-24" after "2025-12-17",
# cloning the layout/styles of the most recent week and replacing
# the ranking data with the new week's values.

$wb = $excel.ActiveWorkbook

# Remember whichever sheet was active before this edit so we can restore
# the selection afterwards (adding/copying a sheet otherwise activates it).
$originallyActiveSheet = $wb.ActiveSheet

# The most recent existing weekly sheet ("2025-12-17") is the last tab.
$sourceSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate it (carries over sheetPr/pageMargins/column formats/header
# style/etc.) and place the copy right after the source sheet.
$sourceSheet.Copy($null, $sourceSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2025-12-24"

# rank, title, volume, highlighted(1/0)
$rows = @(
    ,@(1,'ONE PIECE',113,0)
    ,@(2,'信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐&『ざまぁ!』します!',20,0)
    ,@(3,'しょせん他人事ですから ~とある弁護士の本音の仕事~',10,0)
    ,@(4,'薫る花は凛と咲く',21,0)
    ,@(5,'ケンガンオメガ',31,0)
    ,@(6,'魔入りました!入間くん',46,0)
    ,@(7,'勇者パーティを追い出された器用貧乏 ~パーティ事情で付与術士をやっていた剣士、万能へと至る~',17,0)
    ,@(8,'ワールドトリガー',29,0)
    ,@(9,'ひかえめに言っても、これは愛',8,0)
    ,@(10,'失格紋の最強賢者 ~世界最強の賢者が更に強くなるために転生しました~',33,0)
    ,@(11,'WIND BREAKER',24,0)
    ,@(12,'桃源暗鬼',27,0)
    ,@(13,'誰かこの状況を説明してください! ~契約から始まるウェディング~ 11(アリアンローズコミックス)',11,0)
    ,@(14,'薬屋のひとりごと',16,0)
    ,@(15,'隣のステラ',9,0)
    ,@(16,'ザ・ファブル The third secret',3,1)
    ,@(17,'味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',18,0)
    ,@(18,'SAKAMOTO DAYS',25,0)
    ,@(19,'殲滅魔導の最強賢者 無才の賢者、魔導を極め最強へ至る',10,0)
    ,@(20,'神様のバレー',39,0)
    ,@(21,'最弱テイマーはゴミ拾いの旅を始めました。@COMIC',8,0)
    ,@(22,'じゃあ、あんたが作ってみろよ',1,1)
    ,@(23,'好きじゃないけど、抱いてください1',1,1)
    ,@(24,'推しが武道館いってくれたら死ぬ',12,0)
    ,@(25,'異世界魔王と召喚少女の奴隷魔術',28,0)
    ,@(26,'スライム倒して300年、知らないうちにレベルMAXになってました',17,0)
    ,@(27,'異世界転生の冒険者',13,0)
    ,@(28,'アイツノカノジョ',8,0)
    ,@(29,'FAIRY TAIL 100 YEARS QUEST',22,0)
    ,@(30,'恋愛検証~相性0%男と結婚はアリ?1',1,1)
    ,@(31,'レアモンスター?それ、ただの害虫ですよ',1,1)
    ,@(32,'山、買いました ~異世界暮らしも悪くない~',1,1)
    ,@(33,'組長娘と世話係15',15,0)
    ,@(34,'無能なナナ',14,0)
    ,@(35,'元最強探索者のおじさん。美少女配信者を助けて大バズりしてしまった。',1,1)
    ,@(36,'悪役がいっぱい出てくるエロゲのキモデブ悪役貴族に転生した。痩せて、破滅回避し悪役達による犯罪を未然に防いでスローライフを目指す',1,1)
    ,@(37,'ペストが明けたら遊びましょう! ~中世ヨーロッパ世界と現代文明スローライフ~',1,1)
    ,@(38,'ゆかいな神統記',1,1)
    ,@(39,'近畿地方のある場所について',1,1)
    ,@(40,'悪役令嬢と極道P 異世界のヤクザ、乙女ゲームの悪役令嬢をプロデュースする。',1,1)
    ,@(41,'変な家:',6,0)
    ,@(42,'バーサス',6,0)
    ,@(43,'最強陰陽師の異世界転生記~下僕の妖怪どもに比べてモンスターが弱すぎるんだが~(コミック)',11,0)
    ,@(44,'山口くんはワルくない',12,0)
    ,@(45,'僕らの好きはわりきれない',7,0)
    ,@(46,'怪物事変',24,0)
    ,@(47,'元最強探索者のおじさん。美少女配信者を助けて大バズりしてしまった。',2,1)
    ,@(48,'元最強探索者のおじさん。美少女配信者を助けて大バズりしてしまった。',3,1)
    ,@(49,'悪役がいっぱい出てくるエロゲのキモデブ悪役貴族に転生した。痩せて、破滅回避し悪役達による犯罪を未然に防いでスローライフを目指す',2,1)
    ,@(50,'悪役がいっぱい出てくるエロゲのキモデブ悪役貴族に転生した。痩せて、破滅回避し悪役達による犯罪を未然に防いでスローライフを目指す',3,1)
    ,@(51,'レアモンスター?それ、ただの害虫ですよ',2,1)
    ,@(52,'レアモンスター?それ、ただの害虫ですよ',3,1)
    ,@(53,'近畿地方のある場所について',2,1)
    ,@(54,'近畿地方のある場所について',3,1)
    ,@(55,'山、買いました ~異世界暮らしも悪くない~',2,1)
    ,@(56,'山、買いました ~異世界暮らしも悪くない~',3,1)
    ,@(57,'全員記憶喪失オフィス',1,1)
    ,@(58,'ぷくちょらりファミリア',1,1)
    ,@(59,'じゃあ、あんたが作ってみろよ',4,0)
    ,@(60,'GIANT KILLING',68,0)
    ,@(61,'ちひろさん',10,0)
    ,@(62,'実は俺、最強でした?',18,0)
    ,@(63,'異世界ウォーキング',12,0)
    ,@(64,'ブルーロック',36,0)
    ,@(65,'メガネ、時々、ヤンキーくん',10,0)
    ,@(66,'真綿の檻',7,0)
    ,@(67,'喫茶小鳥の飛び立ちごはん',1,1)
    ,@(68,'すれち恋',3,1)
    ,@(69,'妖魔狩りの末裔-俺だけ不死身の覚醒者-2',2,1)
    ,@(70,'あなたの戸籍、俺にください。4',4,0)
    ,@(71,'ペストが明けたら遊びましょう! ~中世ヨーロッパ世界と現代文明スローライフ~',2,1)
    ,@(72,'ペストが明けたら遊びましょう! ~中世ヨーロッパ世界と現代文明スローライフ~',3,1)
    ,@(73,'ゆかいな神統記',2,1)
    ,@(74,'ゆかいな神統記',3,1)
    ,@(75,'篠原君ちのおうちごはん!',1,1)
    ,@(76,'嶋田と和泉',1,1)
    ,@(77,'全員記憶喪失オフィス',2,1)
    ,@(78,'全員記憶喪失オフィス',3,1)
    ,@(79,'力石持つ',1,1)
    ,@(80,'このマンガがすごい! comics アマテラスの暗号 第一話',1,1)
    ,@(81,'ギャラ飲み女子とラーメンおじさん',2,1)
    ,@(82,'独身貴族は異世界を謳歌する ~結婚しない男の優雅なおひとりさまライフ~',7,0)
    ,@(83,'素材採取家の異世界旅行記9',9,0)
    ,@(84,'魔術ギルド総帥~生まれ変わって今更やり直す2度目の学院生活~',11,0)
    ,@(85,'ザ・ファブル The third secret',1,1)
    ,@(86,'最強は田舎農家のおっさんでした~最高ランクのドラゴンを駆除した結果、実力が世界にバレました~',3,1)
    ,@(87,'MIX',24,0)
    ,@(88,'青の祓魔師',33,0)
    ,@(89,'SPY×FAMILY',16,0)
    ,@(90,'ワンパンマン',35,0)
    ,@(91,'私たちは傷口に愛を塗る1',1,1)
    ,@(92,'あなたの戸籍、俺にください。1',1,1)
    ,@(93,'人気配信者たちのマネージャーになったら、全員元カノだった 第1話',1,1)
    ,@(94,'狙い撃ち 国税調査官 南原&九野 1話',1,1)
    ,@(95,'篠原君ちのおうちごはん!',2,1)
    ,@(96,'篠原君ちのおうちごはん!',3,1)
    ,@(97,'嶋田と和泉',2,1)
    ,@(98,'嶋田と和泉',3,1)
    ,@(99,'力石持つ',2,1)
    ,@(100,'力石持つ',3,1)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]

    $newSheet.Cells.Item($r, 1).Value = $rowData[0]
    $newSheet.Cells.Item($r, 2).Value = $rowData[1]
    $newSheet.Cells.Item($r, 3).Value = $rowData[2]

    if ($rowData[3] -eq 1) {
        $sourceSheet.Range("C5").Copy()
        $newSheet.Cells.Item($r, 3).PasteSpecial(-4122)
    } else {
        $newSheet.Cells.Item($r, 3).ClearFormats()
    }
}

$excel.CutCopyMode = 0

# Restore the original selection/active sheet.
$originallyActiveSheet.Activate()
[void]$originallyActiveSheet.Range("A1").Select()

Write-Host "Added sheet 2025-12-24 with" $rows.Count "rows"
